$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.956.77"
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = "'1.877.98"
$ws.Range("E3").Value = '  -2.04%  '
$ws.Range("D4").Value = "'0.9980"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = "'241.36"
$ws.Range("E5").Value = '  -5.05%  '
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = "'0.4993"
$ws.Range("E7").Value = '  -3.16%  '
$ws.Range("D8").Value = "'44.56"
$ws.Range("E8").Value = '  -3.23%  '
$ws.Range("E9").Value = '  -1.80%  '
$ws.Range("D10").Value = "'0.06612"
$ws.Range("E10").Value = '  -3.38%  '
$ws.Range("D11").Value = "'1.876.76"
$ws.Range("E11").Value = '  -2.04%  '
$ws.Range("E12").Value = '  -4.35%  '
$ws.Range("D13").Value = "'0.07248"
$ws.Range("E13").Value = '  -1.20%  '
$ws.Range("D14").Value = "'0.6667"
$ws.Range("E14").Value = '  -3.36%  '
$ws.Range("D15").Value = "'86.20"
$ws.Range("E15").Value = '  -1.62%  '
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").Value = "'29.929.23"
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").Value = "'0.000007894"
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("D19").Value = "'0.9978"
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("E20").Value = '  -2.48%  '
$ws.Range("D21").Value = "'2.119.23"
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("D22").Value = "'0.9971"
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = "'4.757"
$ws.Range("E23").Value = '  -2.08%  '
$ws.Range("D24").Value = "'5.629"
$ws.Range("E24").Value = '  -2.38%  '
$ws.Range("D25").Value = "'9.070"
$ws.Range("E25").Value = '  -1.32%  '
$ws.Range("D26").Value = "'148.83"
$ws.Range("E26").Value = '  +1.61%  '
$ws.Range("D27").Value = "'138.97"
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("D28").Value = "'16.97"
$ws.Range("E28").Value = '  -1.83%  '
$ws.Range("D29").Value = "'1.908"
$ws.Range("E29").Value = '  -5.27%  '
$ws.Range("D30").Value = "'1.385"
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("D32").Value = "'0.08797"
$ws.Range("E32").Value = '  -0.63%  '
$ws.Range("D33").Value = "'3.951"
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("D34").Value = "'0.05053"
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").Value = "'0.7117"
$ws.Range("E35").Value = '  -1.09%  '
$ws.Range("E36").Value = '  -4.51%  '
$ws.Range("D37").Value = "'2.664"
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("D38").Value = "'2.696"
$ws.Range("E38").Value = '  -4.94%  '
$ws.Range("D39").Value = "'0.01745"
$ws.Range("E39").Value = '  +2.75%  '
$ws.Range("D40").Value = "'2.180"
$ws.Range("E40").Value = '  -6.34%  '
$ws.Range("E41").Value = '  -4.84%  '
$ws.Range("D42").Value = "'0.4262"
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("E43").Value = '  -5.75%  '
$ws.Range("D44").Value = "'0.9971"
$ws.Range("D45").Value = "'101.60"
$ws.Range("E45").Value = '  -3.84%  '
$ws.Range("D46").Value = "'7.435"
$ws.Range("E46").Value = '  -3.49%  '
$ws.Range("D47").Value = "'0.1255"
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("D49").Value = "'32.33"
$ws.Range("E49").Value = '  -3.23%  '
$ws.Range("E50").Value = '  -2.60%  '
$ws.Range("D51").Value = "'8.185"
$ws.Range("E51").Value = '  -4.40%  '
